$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new header row above the existing data (shifts rows 1-6 down to 2-7)
$ws.Range("A1:D1").EntireRow.Insert()

# 2. Populate the new header row
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "CPF"
$ws.Range("C1").Value = "Inscrição"

# 3. Update the "link" column (D) for the first three existing rows with new URLs
$ws.Range("D2").Value = "http://www.africau.edu/images/default/sample.pdf"
$ws.Range("D3").Value = "https://www.w3.org/WAI/ER/tests/xhtml/testfiles/resources/pdf/dummy.pdf"
$ws.Range("D4").Value = "https://juventudedesporto.cplp.org/files/sample-pdf_9359.pdf"

# 4. Replace the candidate in row 5 entirely with new data
$ws.Range("A5").Value = "Ana Clara"
$ws.Range("B5").Value = "'545454"
$ws.Range("C5").Value = 202069
$ws.Range("D5").Value = "https://via.placeholder.com/300.png"

# 5. Clear out the two trailing candidate rows, leaving their formatting intact
$ws.Range("A6:D6").ClearContents()
$ws.Range("A7:D7").ClearContents()

# 6. Mark D9 with an underlined font (placeholder for a future download link), no value
$ws.Range("D9").Font.Underline = $true

# 7. Update the stale AutoFilter defined name range to match the new data layout
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Planilha1!_FilterDatabase") {
        $n.RefersTo = "=Planilha1!`$A`$2:`$D`$7"
    }
}

# 8. Leave the active selection on D9, matching the final saved view
$ws.Range("D9").Select() | Out-Null
